$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value = '43.385.63'
$ws.Range('E2').Value = '  +1.57%  '

$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value = '2.237.60'
$ws.Range('E3').Value = '  +0.58%  '

$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '

$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.30'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.31%  '

$ws.Range('B6').Value = 'Solana'
$ws.Range('C6').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.24'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.54%  '

$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.582'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.73%  '

$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.06%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.562'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.46%  '

$ws.Range('B10').Value = 'Avalanche'
$ws.Range('C10').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.10'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.55%  '

$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0830'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.40%  '

$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.70'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.29%  '

$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.107'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.92%  '

$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.866'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.81%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.37'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.71%  '

$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '2.243.61'
$ws.Range('E16').Value = '  +0.71%  '

$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '43.339.67'
$ws.Range('E17').Value = '  +1.79%  '

$ws.Range('B18').Value = 'InternetComputer(DFINITY)'
$ws.Range('C18').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '14.16'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.86%  '

$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.65'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.41%  '

$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0972'
$ws.Range('E20').Value = '  +2.85%  '

$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '65.30'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.77%  '

$ws.Range('B22').Value = 'PancakeSwap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.19'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.43%  '

$ws.Range('B23').Value = 'BitcoinCash'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.29'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.93%  '

$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.18'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.63%  '

$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.14%  '

$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.05'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.10%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.01%  '

$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.22'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.13%  '

$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '36.88'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.13%  '

$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.39'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.30%  '

$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.28'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.78%  '

$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0874'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.33%  '

$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '157.35'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.40%  '

$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.71'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.00%  '

$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.27'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +4.83%  '

$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.121'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.52%  '

$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.89'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.91%  '

$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.41'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.39%  '

$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.104'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.57%  '

$ws.Range('B40').Value = 'NEARProtocol'
$ws.Range('C40').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.70'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.15%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0322'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.47%  '

$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '14.38'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +18.37%  '

$ws.Range('B43').Value = 'FirstDigitalUSD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.02%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.818.75'
$ws.Range('E44').Value = '  +0.55%  '

$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.203'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.55%  '

$ws.Range('B46').Value = 'BitcoinSV'
$ws.Range('C46').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '84.13'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -5.81%  '

$ws.Range('B47').Value = 'THORChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.31'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.63%  '

$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.86'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.28%  '

$ws.Range('B49').Value = 'ordi'
$ws.Range('C49').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '73.95'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.83%  '

$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '103.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.63%  '

$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '58.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.49%  '

